$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cell C10 value changed from 18 to 1 (numeric)
$ws.Range("C10").Value = 1
